$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.62799999999998
$ws.Range("A7").Value = -19.33589999999999
$ws.Range("C7").Value = -12.5133
$ws.Range("C15").Value = -14.26
$ws.Range("A16").Value = -21.53169999999998
$ws.Range("D16").Value = -9.160700000000009
$ws.Range("D19").Value = -8.827899999999996
$ws.Range("C21").Value = -11.96890000000001
$ws.Range("C22").Value = -12.33240000000001
$ws.Range("C23").Value = -12.64890000000001
$ws.Range("A28").Value = -21.97199999999999
$ws.Range("A29").Value = -21.11149999999997
$ws.Range("A32").Value = -21.26679999999999
$ws.Range("C34").Value = -11.74620000000001
$ws.Range("E34").Value = 17.5533
$ws.Range("D36").Value = -8.720399999999994
$ws.Range("A40").Value = -20.4202
$ws.Range("C43").Value = -12.31199999999999
$ws.Range("E43").Value = 17.26750000000002
$ws.Range("C45").Value = -14.02469999999999
$ws.Range("D46").Value = -8.649900000000001
$ws.Range("E48").Value = 17.4418
$ws.Range("C50").Value = -14.18169999999999
$ws.Range("D50").Value = -7.997899999999996
$ws.Range("C51").Value = -11.9792
$ws.Range("A52").Value = -22.1478
$ws.Range("A57").Value = -22.41760000000001
$ws.Range("A66").Value = -22.0738
$ws.Range("C66").Value = -12.9632
$ws.Range("C67").Value = -11.055
$ws.Range("E70").Value = 17.61820000000001
$ws.Range("E73").Value = 17.43660000000001
$ws.Range("C79").Value = -12.0569
$ws.Range("C84").Value = -13.43299999999999
$ws.Range("E87").Value = 16.31409999999999
$ws.Range("C92").Value = -11.6854
$ws.Range("E92").Value = 18.30950000000001
$ws.Range("D95").Value = -8.329999999999991
$ws.Range("C97").Value = -11.496
$ws.Range("D97").Value = -8.449299999999994
$ws.Range("A100").Value = -22.0565
$ws.Range("E101").Value = 16.69620000000001
